$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.79399266666667
$ws.Range("H2").Value = 47.381978
$ws.Range("I2").Value = 0.2968109173698557
$ws.Range("J2").Value = 0.2968109173698557
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.564139666666667
$ws.Range("N2").Value = 4.692419
$ws.Range("O2").Value = 0.166125338305886
$ws.Range("P2").Value = 0.166125338305886
$ws.Range("Q2").Value = 24.70401042497578
$ws.Range("R2").Value = 222.336093824782
$ws.Range("S2").Value = 0.04930781406094766
$ws.Range("T2").Value = 0.04930781406094764

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.79399266666667
$ws.Range("H3").Value = 47.381978
$ws.Range("I3").Value = 0.2968109173698557
$ws.Range("J3").Value = 0.2968109173698557
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.316850333333333
$ws.Range("N3").Value = 9.950551
$ws.Range("O3").Value = 0.3522785691569683
$ws.Range("P3").Value = 0.3522785691569683
$ws.Range("Q3").Value = 52.38630984109756
$ws.Range("R3").Value = 471.4767885698781
$ws.Range("S3").Value = 0.1045601252812199
$ws.Range("T3").Value = 0.1045601252812199

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.79399266666667
$ws.Range("H4").Value = 47.381978
$ws.Range("I4").Value = 0.2968109173698557
$ws.Range("J4").Value = 0.2968109173698557
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.534428999999999
$ws.Range("N4").Value = 13.603287
$ws.Range("O4").Value = 0.4815960925371456
$ws.Range("P4").Value = 0.4815960925371456
$ws.Range("Q4").Value = 71.61673837352066
$ws.Range("R4").Value = 644.550645361686
$ws.Range("S4").Value = 0.1429429780276881
$ws.Range("T4").Value = 0.1429429780276881

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 24.86954866666666
$ws.Range("H5").Value = 74.608646
$ws.Range("I5").Value = 0.4673646309781075
$ws.Range("J5").Value = 0.4673646309781075
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.564139666666667
$ws.Range("N5").Value = 4.692419
$ws.Range("O5").Value = 0.166125338305886
$ws.Range("P5").Value = 0.166125338305886
$ws.Range("Q5").Value = 38.89944756163044
$ws.Range("R5").Value = 350.095028054674
$ws.Range("S5").Value = 0.07764110743344368
$ws.Range("T5").Value = 0.07764110743344368

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 24.86954866666666
$ws.Range("H6").Value = 74.608646
$ws.Range("I6").Value = 0.4673646309781075
$ws.Range("J6").Value = 0.4673646309781075
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.316850333333333
$ws.Range("N6").Value = 9.950551
$ws.Range("O6").Value = 0.3522785691569683
$ws.Range("P6").Value = 0.3522785691569683
$ws.Range("Q6").Value = 82.48857078488288
$ws.Range("R6").Value = 742.397137063946
$ws.Range("S6").Value = 0.1646425434755422
$ws.Range("T6").Value = 0.1646425434755422

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.86954866666666
$ws.Range("H7").Value = 74.608646
$ws.Range("I7").Value = 0.4673646309781075
$ws.Range("J7").Value = 0.4673646309781075
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.534428999999999
$ws.Range("N7").Value = 13.603287
$ws.Range("O7").Value = 0.4815960925371456
$ws.Range("P7").Value = 0.4815960925371456
$ws.Range("Q7").Value = 112.7692026910446
$ws.Range("R7").Value = 1014.922824219402
$ws.Range("S7").Value = 0.2250809800691216
$ws.Range("T7").Value = 0.2250809800691216

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.54876233333333
$ws.Range("H8").Value = 37.646287
$ws.Range("I8").Value = 0.2358244516520368
$ws.Range("J8").Value = 0.2358244516520368
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.564139666666667
$ws.Range("N8").Value = 4.692419
$ws.Range("O8").Value = 0.166125338305886
$ws.Range("P8").Value = 0.166125338305886
$ws.Range("Q8").Value = 19.62801693313923
$ws.Range("R8").Value = 176.652152398253
$ws.Range("S8").Value = 0.03917641681149468
$ws.Range("T8").Value = 0.03917641681149467

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.54876233333333
$ws.Range("H9").Value = 37.646287
$ws.Range("I9").Value = 0.2358244516520368
$ws.Range("J9").Value = 0.2358244516520368
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.316850333333333
$ws.Range("N9").Value = 9.950551
$ws.Range("O9").Value = 0.3522785691569683
$ws.Range("P9").Value = 0.3522785691569683
$ws.Range("Q9").Value = 41.62236652823745
$ws.Range("R9").Value = 374.601298754137
$ws.Range("S9").Value = 0.08307590040020618
$ws.Range("T9").Value = 0.08307590040020618

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.54876233333333
$ws.Range("H10").Value = 37.646287
$ws.Range("I10").Value = 0.2358244516520368
$ws.Range("J10").Value = 0.2358244516520368
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.534428999999999
$ws.Range("N10").Value = 13.603287
$ws.Range("O10").Value = 0.4815960925371456
$ws.Range("P10").Value = 0.4815960925371456
$ws.Range("Q10").Value = 56.90147183837433
$ws.Range("R10").Value = 512.113246545369
$ws.Range("S10").Value = 0.113572134440336
$ws.Range("T10").Value = 0.1135721344403359

